$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2; D = "66.839.83"; E = "  -3.96%  " },
    @{ Row = 3; D = "3.347.35"; E = "  -0.79%  " },
    @{ Row = 4; E = "  +0.01%  " },
    @{ Row = 5; D = "573.24"; E = "  -3.49%  " },
    @{ Row = 6; D = "181.94"; E = "  -5.19%  " },
    @{ Row = 7; E = "  +0.04%  " },
    @{ Row = 8; D = "0.599"; E = "  -1.60%  " },
    @{ Row = 9; E = "  -3.79%  " },
    @{ Row = 10; D = "6.65" },
    @{ Row = 12; D = "3.930.37"; E = "  -0.79%  " },
    @{ Row = 13; E = "  -1.60%  " },
    @{ Row = 14; D = "27.18"; E = "  -5.35%  " },
    @{ Row = 15; D = "66.885.13"; E = "  -3.92%  " },
    @{ Row = 16; E = "  -2.56%  " },
    @{ Row = 17; D = "3.341.16"; E = "  -1.08%  " },
    @{ Row = 18; D = "436.28"; E = "  -3.19%  " },
    @{ Row = 19; D = "13.65"; E = "  -1.58%  " },
    @{ Row = 20; E = "  -2.76%  " },
    @{ Row = 21; D = "7.64"; E = "  -2.52%  " },
    @{ Row = 22; D = "73.86"; E = "  -1.24%  " },
    @{ Row = 23; D = "0.999"; E = "  -0.01%  " },
    @{ Row = 24; D = "0.520"; E = "  +0.08%  " },
    @{ Row = 25; E = "  -3.83%  " },
    @{ Row = 26; E = "  -0.55%  " },
    @{ Row = 27; D = "9.09"; E = "  -5.22%  " },
    @{ Row = 28; E = "  +0.29%  " },
    @{ Row = 29; D = "1.96"; E = "  -1.96%  " },
    @{ Row = 30; D = "22.92"; E = "  -1.71%  " },
    @{ Row = 31; D = "5.31"; E = "  -6.30%  " },
    @{ Row = 32; E = "  +0.00%  " },
    @{ Row = 33; E = "  -4.31%  " },
    @{ Row = 34; D = "6.80"; E = "  -3.41%  " },
    @{ Row = 35; D = "1.50"; E = "  -2.21%  " },
    @{ Row = 36; D = "161.24"; E = "  -2.40%  " },
    @{ Row = 37; D = "27.78"; E = "  +1.87%  " },
    @{ Row = 38; D = "1.85"; E = "  -4.93%  " },
    @{ Row = 39; D = "2.845.05"; E = "  +3.86%  " },
    @{ Row = 40; D = "0.794"; E = "  -3.11%  " },
    @{ Row = 41; D = "4.44"; E = "  -4.22%  " },
    @{ Row = 42; D = "6.23"; E = "  -5.08%  " },
    @{ Row = 43; D = "0.0674"; E = "  -2.66%  " },
    @{ Row = 44; D = "40.24"; E = "  -1.06%  " },
    @{ Row = 45; D = "24.61"; E = "  -4.63%  " },
    @{ Row = 46; D = "2.37"; E = "  -7.22%  " },
    @{ Row = 47; D = "324.52"; E = "  -5.72%  " },
    @{ Row = 48; E = "  -4.69%  " },
    @{ Row = 49; D = "31.38"; E = "  -5.25%  " },
    @{ Row = 50; D = "0.989"; E = "  -3.90%  " },
    @{ Row = 51; D = "6.16"; E = "  -2.90%  " }
)

foreach ($item in $data) {
    if ($item.ContainsKey("D")) {
        $cell = $ws.Range("D" + $item.Row)
        $cell.NumberFormat = "@"
        $cell.Value = $item.D
        $cell.Style = "Normal"
    }
    if ($item.ContainsKey("E")) {
        $ws.Range("E" + $item.Row).Value = $item.E
    }
}
